# Auto-generated edit script applying the cryptos.xlsx data refresh
# (GitHub Actions scheduled update) described by the target diff.
# Cells whose new text is a plain numeric token (e.g. "605.47") are
# written with a leading apostrophe so Excel keeps them as literal
# text (matching the inlineStr/shared-string cells in the workbook)
# instead of silently re-typing them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.425.26"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "3.555.55"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'605.47"
$ws.Range("D6").Value = "'144.76"
$ws.Range("E6").Value = "  +0.80%  "
$ws.Range("D7").Value = "3.553.34"
$ws.Range("E7").Value = "  +0.57%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").Value = "'0.497"
$ws.Range("E9").Value = "  +3.46%  "
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("D11").Value = "'7.94"
$ws.Range("E11").Value = "  -1.33%  "
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").Value = "4.156.12"
$ws.Range("E13").Value = "  +0.64%  "
$ws.Range("E14").Value = "  -0.39%  "
$ws.Range("D15").Value = "'29.99"
$ws.Range("E15").Value = "  -0.66%  "
$ws.Range("D16").Value = "3.591.88"
$ws.Range("E16").Value = "  +1.77%  "
$ws.Range("D17").Value = "66.470.33"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").Value = "'11.59"
$ws.Range("E19").Value = "  +5.72%  "
$ws.Range("E20").Value = "  -0.60%  "
$ws.Range("E21").Value = "  -0.59%  "
$ws.Range("D22").Value = "'431.20"
$ws.Range("E22").Value = "  +1.32%  "
$ws.Range("D23").Value = "'0.610"
$ws.Range("E23").Value = "  +1.33%  "
$ws.Range("D24").Value = "'79.66"
$ws.Range("E24").Value = "  +1.30%  "
$ws.Range("D25").Value = "3.696.36"
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").Value = "  -0.66%  "
$ws.Range("E28").Value = "  +0.54%  "
$ws.Range("E29").Value = "  +1.09%  "
$ws.Range("E30").Value = "  -1.26%  "
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("D32").Value = "'1.46"
$ws.Range("E32").Value = "  -2.22%  "
$ws.Range("D33").Value = "3.548.32"
$ws.Range("E33").Value = "  +0.76%  "
$ws.Range("D34").Value = "'25.34"
$ws.Range("E34").Value = "  +0.32%  "
$ws.Range("E35").Value = "  -4.88%  "
$ws.Range("E36").Value = "  +0.60%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("E38").Value = "  -2.07%  "
$ws.Range("E39").Value = "  -0.86%  "
$ws.Range("D40").Value = "'174.80"
$ws.Range("E40").Value = "  +1.66%  "
$ws.Range("D41").Value = "'0.0847"
$ws.Range("E41").Value = "  -1.14%  "
$ws.Range("D42").Value = "'5.19"
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("E43").Value = "  -0.46%  "
$ws.Range("E44").Value = "  +1.44%  "
$ws.Range("D45").Value = "'46.14"
$ws.Range("E45").Value = "  +1.48%  "
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").Value = "'2.52"
$ws.Range("E47").Value = "  +4.93%  "
$ws.Range("E48").Value = "  -2.34%  "
$ws.Range("D49").Value = "'25.10"
$ws.Range("E49").Value = "  -3.62%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'23.51"
$ws.Range("E50").Value = "  +4.77%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "'7.14"
$ws.Range("E51").Value = "  +0.03%  "
